$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the standalone row that only held the "5840535 - Messias Borges Silva"
# value (it used to sit directly under "Docentes responsáveis:"). Deleting it
# shifts every following row up by one.
$ws.Rows(13).Delete()

# After the shift, re-point the remaining value rows to their new text.
$ws.Range("B10").Value = "5840535 - Messias Borges Silva"
$ws.Range("C10").Value = "5840535 - Messias Borges Silva"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# "01/01/2022" must stay plain text (it already is, two rows above in B8/C8).
# Assigning it through .Value would make Excel reinterpret it as a date
# serial, so copy the already-text cell instead of retyping the literal.
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))

$ws.Range("B18").Value = "5840535 - Messias Borges Silva"
$ws.Range("C18").Value = "5840535 - Messias Borges Silva"

$ws.Range("B19").Value = "Provas, relatórios e apresentação de seminários."
$ws.Range("C19").Value = "Provas, relatórios e apresentação de seminários."

$ws.Range("B20").Value = "MF = (0,7*P&R + 0,3*S), onde P&R= Prova e relatórios e S= Seminário."
$ws.Range("C20").Value = "MF = (0,7*P&R + 0,3*S), onde P&R= Prova e relatórios e S= Seminário."

$ws.Range("B21").Value = "É feita sob a forma de uma prova, com toda a matéria dada, com duas horas de duração, aplicada no período determinado pela USP. A média final será a média aritmética entre a nota desta prova e a média obtida no semestre."
$ws.Range("C21").Value = "É feita sob a forma de uma prova, com toda a matéria dada, com duas horas de duração, aplicada no período determinado pela USP. A média final será a média aritmética entre a nota desta prova e a média obtida no semestre."
